# Generate Report for Handoff
# b.md has been handed off for zh-cn and de-de: update the status/report
# rows (row 3 on each sheet) to reflect the new handoff file + timestamps.
$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (b.md) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-29-18 20:29:26"

# --- zh-cn sheet: row 3 (b.md) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-18 20:29:24"

# --- de-de sheet: row 3 (b.md) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-18 20:29:26"
